$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9999999967220834,
    0.99479307412868612,
    0.97726943032200508,
    0.97018939803775406,
    0.9636065023774224,
    0.94792193459981955,
    0.94570514729175459,
    0.94392585551094954,
    0.94069633902744321,
    0.93840589649286121,
    0.93821143121777073,
    0.93809764889493419,
    0.93901401801028883,
    0.94009143001250028,
    0.94254980926845222,
    0.94004351624336246,
    0.93633576309574162,
    0.93522684407450118,
    0.99373981757384477,
    0.98662262504882681,
    0.98522411254767506,
    0.98395960341563737,
    0.97023406161970516,
    0.95721261279194358,
    0.95075551494281152,
    0.94599558504467929,
    0.94356004578486685,
    0.93345096134749606,
    0.92662784680160626,
    0.92420847579428833,
    0.9273608447751307,
    0.92568151902595863,
    0.92516149888589572
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
